$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last charge end time" column (D) for all data rows (2-54) to the new refresh timestamp
$newRefreshTime = 45983.327719907407
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 4).Value = $newRefreshTime
}

# Rows 18-54: station name (A), terminal name (B) and last-charge-start-time (C) are refreshed
# with the latest outstanding list of terminals that have not charged for a long time.
$ws.Cells.Item(18, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(18, 2).Value = "101号直流"
$ws.Cells.Item(18, 3).Value = 45979.18986111111
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "603号直流"
$ws.Cells.Item(19, 3).Value = 45980.25017361111
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "904号直流"
$ws.Cells.Item(20, 3).Value = 45981.57524305556
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "B01号直流"
$ws.Cells.Item(21, 3).Value = 45981.67456018519
$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value = "902号直流"
$ws.Cells.Item(22, 3).Value = 45981.738391203704
$ws.Cells.Item(23, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(23, 2).Value = "110号直流"
$ws.Cells.Item(23, 3).Value = 45982.02722222222
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "105号直流"
$ws.Cells.Item(24, 3).Value = 45982.043587962966
$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(25, 2).Value = "B02号直流"
$ws.Cells.Item(25, 3).Value = 45982.09847222222
$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(26, 2).Value = "306号直流"
$ws.Cells.Item(26, 3).Value = 45982.115694444445
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(27, 2).Value = "106号直流"
$ws.Cells.Item(27, 3).Value = 45982.1621875
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value = "103号直流"
$ws.Cells.Item(28, 3).Value = 45982.25917824074
$ws.Cells.Item(29, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(29, 2).Value = "207号直流"
$ws.Cells.Item(29, 3).Value = 45982.51305555556
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(30, 2).Value = "203号直流"
$ws.Cells.Item(30, 3).Value = 45982.53803240741
$ws.Cells.Item(31, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "805号直流"
$ws.Cells.Item(31, 3).Value = 45982.54555555555
$ws.Cells.Item(32, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(32, 2).Value = "505号直流"
$ws.Cells.Item(32, 3).Value = 45982.55150462963
$ws.Cells.Item(33, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(33, 2).Value = "402号直流"
$ws.Cells.Item(33, 3).Value = 45982.553715277776
$ws.Cells.Item(34, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(34, 2).Value = "502号直流"
$ws.Cells.Item(34, 3).Value = 45982.55546296296
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value = "002A号直流"
$ws.Cells.Item(35, 3).Value = 45982.555972222224
$ws.Cells.Item(36, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(36, 2).Value = "A04号直流"
$ws.Cells.Item(36, 3).Value = 45982.55740740741
$ws.Cells.Item(37, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(37, 2).Value = "004A号直流"
$ws.Cells.Item(37, 3).Value = 45982.55820601852
$ws.Cells.Item(38, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(38, 2).Value = "602号直流"
$ws.Cells.Item(38, 3).Value = 45982.558599537035
$ws.Cells.Item(39, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(39, 2).Value = "301号直流"
$ws.Cells.Item(39, 3).Value = 45982.573229166665
$ws.Cells.Item(40, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(40, 2).Value = "903号直流"
$ws.Cells.Item(40, 3).Value = 45982.57377314815
$ws.Cells.Item(41, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(41, 2).Value = "406号直流"
$ws.Cells.Item(41, 3).Value = 45982.584861111114
$ws.Cells.Item(42, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(42, 2).Value = "103号直流"
$ws.Cells.Item(42, 3).Value = 45982.61107638889
$ws.Cells.Item(43, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(43, 2).Value = "803号直流"
$ws.Cells.Item(43, 3).Value = 45982.61452546297
$ws.Cells.Item(44, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(44, 2).Value = "604号直流"
$ws.Cells.Item(44, 3).Value = 45982.639398148145
$ws.Cells.Item(45, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(45, 2).Value = "B03号直流"
$ws.Cells.Item(45, 3).Value = 45982.672372685185
$ws.Cells.Item(46, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(46, 2).Value = "401号直流"
$ws.Cells.Item(46, 3).Value = 45982.714004629626
$ws.Cells.Item(47, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(47, 2).Value = "403号直流"
$ws.Cells.Item(47, 3).Value = 45982.719826388886
$ws.Cells.Item(48, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(48, 2).Value = "802号直流"
$ws.Cells.Item(48, 3).Value = 45982.733668981484
$ws.Cells.Item(49, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(49, 2).Value = "A01号直流"
$ws.Cells.Item(49, 3).Value = 45982.73400462963
$ws.Cells.Item(50, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(50, 2).Value = "108号直流"
$ws.Cells.Item(50, 3).Value = 45982.7346875
$ws.Cells.Item(51, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(51, 2).Value = "104号直流"
$ws.Cells.Item(51, 3).Value = 45982.73706018519
$ws.Cells.Item(52, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(52, 2).Value = "211号直流"
$ws.Cells.Item(52, 3).Value = 45982.771099537036
$ws.Cells.Item(53, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(53, 2).Value = "901号直流"
$ws.Cells.Item(53, 3).Value = 45982.77679398148
$ws.Cells.Item(54, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(54, 2).Value = "304号直流"
$ws.Cells.Item(54, 3).Value = 45982.79420138889

# Move the active selection from E12 to E9 as recorded for this save
[void]$ws.Range("E9").Select()
